$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021596613382186
$ws.Range("D2").Value = 1.029479468168048
$ws.Range("E2").Value = 1.022447485649909
$ws.Range("F2").Value = 1.020030238003467
$ws.Range("I2").Value = 1.031460032946297
$ws.Range("J2").Value = 1.026787187609181
$ws.Range("K2").Value = 1.032293259289666
$ws.Range("L2").Value = 1.025281839128513
$ws.Range("M2").Value = 1.022871731065522
$ws.Range("N2").Value = 1.028245343189724
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022572303617348
$ws.Range("D3").Value = 1.029943028407246
$ws.Range("E3").Value = 1.023276204496457
$ws.Range("F3").Value = 1.021644933153752
$ws.Range("I3").Value = 1.031611535565965
$ws.Range("J3").Value = 1.027400428603253
$ws.Range("K3").Value = 1.03256590163614
$ws.Range("L3").Value = 1.02591716585047
$ws.Range("M3").Value = 1.024290359326969
$ws.Range("N3").Value = 1.028859455056346
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023203531793888
$ws.Range("D4").Value = 1.030243103508403
$ws.Range("E4").Value = 1.023812723748472
$ws.Range("F4").Value = 1.022689537233725
$ws.Range("I4").Value = 1.031708358163379
$ws.Range("J4").Value = 1.027796555992132
$ws.Range("K4").Value = 1.03274171451911
$ws.Range("L4").Value = 1.026327904619073
$ws.Range("M4").Value = 1.025207635412995
$ws.Range("N4").Value = 1.029256144991554
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023468874744357
$ws.Range("D5").Value = 1.030369282598172
$ws.Range("E5").Value = 1.024038344088202
$ws.Range("F5").Value = 1.023128644410596
$ws.Range("I5").Value = 1.031748772546693
$ws.Range("J5").Value = 1.027962925412578
$ws.Range("K5").Value = 1.032815480706571
$ws.Range("L5").Value = 1.026500493171662
$ws.Range("M5").Value = 1.025593102989395
$ws.Range("N5").Value = 1.029422750675659
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02351342552372
$ws.Range("D6").Value = 1.030390470196415
$ws.Range("E6").Value = 1.02407623070145
$ws.Range("F6").Value = 1.023202370017978
$ws.Range("I6").Value = 1.031755541296742
$ws.Range("J6").Value = 1.027990850050481
$ws.Range("K6").Value = 1.032827857817872
$ws.Range("L6").Value = 1.02652946650972
$ws.Range("M6").Value = 1.02565781574233
$ws.Range("N6").Value = 1.029450714969751
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023207077418299
$ws.Range("D7").Value = 1.030244789413186
$ws.Range("E7").Value = 1.023815738232763
$ws.Range("F7").Value = 1.022695404771876
$ws.Range("I7").Value = 1.031708899321492
$ws.Range("J7").Value = 1.027798779666648
$ws.Range("K7").Value = 1.032742700759343
$ws.Range("L7").Value = 1.026330211093568
$ws.Range("M7").Value = 1.025212786653089
$ws.Range("N7").Value = 1.029258371823942
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021926375345831
$ws.Range("D8").Value = 1.029636104728316
$ws.Range("E8").Value = 1.022727496496734
$ws.Range("F8").Value = 1.020575979404959
$ws.Range("I8").Value = 1.031511484313389
$ws.Range("J8").Value = 1.026994576420683
$ws.Range("K8").Value = 1.032385524872062
$ws.Range("L8").Value = 1.025496625443229
$ws.Range("M8").Value = 1.023351305377662
$ws.Range("N8").Value = 1.028453026517124
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019668753426308
$ws.Range("D9").Value = 1.028564508175388
$ws.Range("E9").Value = 1.020812042057225
$ws.Range("F9").Value = 1.016839378267776
$ws.Range("I9").Value = 1.031154358484228
$ws.Range("J9").Value = 1.025572237392057
$ws.Range("K9").Value = 1.031751533805329
$ws.Range("L9").Value = 1.024024974153781
$ws.Range("M9").Value = 1.020065754045933
$ws.Range("N9").Value = 1.027028667603933
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018163049409088
$ws.Range("D10").Value = 1.027850851831494
$ws.Range("E10").Value = 1.019536527405738
$ws.Range("F10").Value = 1.014346636514582
$ws.Range("I10").Value = 1.030910064015975
$ws.Range("J10").Value = 1.024620467346704
$ws.Range("K10").Value = 1.031325821274979
$ws.Range("L10").Value = 1.023041998053098
$ws.Range("M10").Value = 1.017871426868625
$ws.Range("N10").Value = 1.026075545935946
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017510904588817
$ws.Range("D11").Value = 1.027542023313583
$ws.Range("E11").Value = 1.018984559223466
$ws.Range("F11").Value = 1.013266766523025
$ws.Range("I11").Value = 1.030802811017418
$ws.Range("J11").Value = 1.024207492941855
$ws.Range("K11").Value = 1.031140767831423
$ws.Range("L11").Value = 1.022615909658112
$ws.Range("M11").Value = 1.016920243791751
$ws.Range("N11").Value = 1.025661985060073
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01726864314521
$ws.Range("D12").Value = 1.027427340442024
$ws.Range("E12").Value = 1.018779584238594
$ws.Range("F12").Value = 1.012865572504695
$ws.Range("I12").Value = 1.030762751480388
$ws.Range("J12").Value = 1.02405396722269
$ws.Range("K12").Value = 1.031071923672715
$ws.Range("L12").Value = 1.022457572918018
$ws.Range("M12").Value = 1.016566771449789
$ws.Range("N12").Value = 1.025508241316778
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017320610232215
$ws.Range("D13").Value = 1.027451938937831
$ws.Range("E13").Value = 1.018823549783663
$ws.Range("F13").Value = 1.012951633827449
$ws.Range("I13").Value = 1.030771354385538
$ws.Range("J13").Value = 1.024086904842733
$ws.Range("K13").Value = 1.031086695812963
$ws.Range("L13").Value = 1.022491539794956
$ws.Range("M13").Value = 1.01664259984082
$ws.Range("N13").Value = 1.025541225712019
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017490879715002
$ws.Range("D14").Value = 1.027532542976233
$ws.Range("E14").Value = 1.018967614896869
$ws.Range("F14").Value = 1.013233605422003
$ws.Range("I14").Value = 1.030799504191106
$ws.Range("J14").Value = 1.024194805086146
$ws.Range("K14").Value = 1.031135079332154
$ws.Range("L14").Value = 1.022602822894536
$ws.Range("M14").Value = 1.016891028972161
$ws.Range("N14").Value = 1.025649279186154
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017595784910684
$ws.Range("D15").Value = 1.027582209771906
$ws.Range("E15").Value = 1.019056384862768
$ws.Range("F15").Value = 1.013407326302947
$ws.Range("I15").Value = 1.030816818935021
$ws.Range("J15").Value = 1.024261268921065
$ws.Range("K15").Value = 1.031164875847183
$ws.Range("L15").Value = 1.022671378989599
$ws.Range("M15").Value = 1.017044072880349
$ws.Range("N15").Value = 1.025715837407342
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018206326548456
$ws.Range("D16").Value = 1.027871351867707
$ws.Range("E16").Value = 1.019573166820676
$ws.Range("F16").Value = 1.014418292819574
$ws.Range("I16").Value = 1.030917151045044
$ws.Range("J16").Value = 1.024647857095774
$ws.Range("K16").Value = 1.03133808761101
$ws.Range("L16").Value = 1.023070266549458
$ws.Range("M16").Value = 1.017934531599075
$ws.Range("N16").Value = 1.026102974581602
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.018589258280126
$ws.Range("D17").Value = 1.028052774717431
$ws.Range("E17").Value = 1.019897420852446
$ws.Range("F17").Value = 1.015052307128846
$ws.Range("I17").Value = 1.030979692759064
$ws.Range("J17").Value = 1.024890125144328
$ws.Range("K17").Value = 1.031446547256959
$ws.Range("L17").Value = 1.023320356518564
$ws.Range("M17").Value = 1.018492813769024
$ws.Range("N17").Value = 1.026345586678573
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018812600025284
$ws.Range("D18").Value = 1.028158613780402
$ws.Range("E18").Value = 1.020086585378143
$ws.Range("F18").Value = 1.015422069571934
$ws.Range("I18").Value = 1.031016030258534
$ws.Range("J18").Value = 1.0250313538681
$ws.Range("K18").Value = 1.031509740698463
$ws.Range("L18").Value = 1.023466186063758
$ws.Range("M18").Value = 1.018818351867709
$ws.Range("N18").Value = 1.026487015963333
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018888751130427
$ws.Range("D19").Value = 1.028194705209256
$ws.Range("E19").Value = 1.020151091134882
$ws.Range("F19").Value = 1.015548141205737
$ws.Range("I19").Value = 1.031028396312874
$ws.Range("J19").Value = 1.025079495303101
$ws.Range("K19").Value = 1.031531276265104
$ws.Range("L19").Value = 1.023515902770692
$ws.Range("M19").Value = 1.018929335449291
$ws.Range("N19").Value = 1.026535225764693
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018548174976512
$ws.Range("D20").Value = 1.028033307872057
$ws.Range("E20").Value = 1.019862628089383
$ws.Range("F20").Value = 1.014984288348547
$ws.Range("I20").Value = 1.030972997314243
$ws.Range("J20").Value = 1.024864140575365
$ws.Range("K20").Value = 1.031434917719843
$ws.Range("L20").Value = 1.023293528755111
$ws.Range("M20").Value = 1.018432925605443
$ws.Range("N20").Value = 1.026319565208542
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.01744074031722
$ws.Range("D21").Value = 1.027508806262876
$ws.Range("E21").Value = 1.018925189917459
$ws.Range("F21").Value = 1.013150574146714
$ws.Range("I21").Value = 1.03079122087453
$ws.Range("J21").Value = 1.024163034710331
$ws.Range("K21").Value = 1.031120834542016
$ws.Range("L21").Value = 1.02257005469096
$ws.Range("M21").Value = 1.016817877271399
$ws.Range("N21").Value = 1.025617463692762
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016744302289714
$ws.Range("D22").Value = 1.027179204296788
$ws.Range("E22").Value = 1.018336077824564
$ws.Range("F22").Value = 1.011997165721287
$ws.Range("I22").Value = 1.03067565225158
$ws.Range("J22").Value = 1.023721477026349
$ws.Range("K22").Value = 1.030922738832308
$ws.Range("L22").Value = 1.022114781064001
$ws.Range("M22").Value = 1.015801499828277
$ws.Range("N22").Value = 1.025175278946221
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017113511818647
$ws.Range("D23").Value = 1.027353915692083
$ws.Range("E23").Value = 1.018648349642043
$ws.Range("F23").Value = 1.012608657369169
$ws.Range("I23").Value = 1.030737038513435
$ws.Range("J23").Value = 1.023955625918058
$ws.Range("K23").Value = 1.031027811583276
$ws.Range("L23").Value = 1.022356167986658
$ws.Range("M23").Value = 1.016340391344133
$ws.Range("N23").Value = 1.025409760356214
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018566738801646
$ws.Range("D24").Value = 1.028042104045498
$ws.Range("E24").Value = 1.019878349340357
$ws.Range("F24").Value = 1.015015023252148
$ws.Range("I24").Value = 1.030976023136355
$ws.Range("J24").Value = 1.024875882137098
$ws.Range("K24").Value = 1.031440172821051
$ws.Range("L24").Value = 1.023305651201454
$ws.Range("M24").Value = 1.018459986791335
$ws.Range("N24").Value = 1.02633132344464
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020252509340196
$ws.Range("D25").Value = 1.028841416721748
$ws.Range("E25").Value = 1.021306976334835
$ws.Range("F25").Value = 1.017805645294668
$ws.Range("I25").Value = 1.031247779612939
$ws.Range("J25").Value = 1.025940568734609
$ws.Range("K25").Value = 1.031915976648548
$ws.Range("L25").Value = 1.024405760686095
$ws.Range("M25").Value = 1.020915820721069
$ws.Range("N25").Value = 1.02739752201924
